$wb = $excel.ActiveWorkbook

# Fill in the data cells on the "AddGroup" sheet
$ws = $wb.Worksheets.Item("AddGroup")

$ws.Cells.Item(2, 1).Value = 1111111
$ws.Cells.Item(2, 2).Value = 2014
$ws.Cells.Item(2, 3).Value = 2018

$ws.Cells.Item(3, 1).Value = 2222222
$ws.Cells.Item(3, 2).Value = 2014
$ws.Cells.Item(3, 3).Value = 2018

# Make "AddGroup" the active sheet, with C2 selected
$ws.Activate()
$ws.Range("C2").Select()

